# Atualização dos dados: 24.12.2025 09:02
# Adds a new data row (row 12, id=11) to the "quantidade_pontos" sheet,
# mirroring the layout/formatting of the preceding rows, and moves the
# active-cell selection as recorded in the session.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy per-cell number formats from row 11 (the previous last row) so the
# new row's Tempo (H) and Data (L) cells pick up the same time/date formats
# instead of the plain column defaults. ---
$ws.Range("H11").Copy()
$ws.Range("H12").PasteSpecial(-4122)

$ws.Range("L11").Copy()
$ws.Range("L12").PasteSpecial(-4122)

# --- Fill in the new row's values ---
$ws.Range("B12").Value = 15
$ws.Range("C12").Value = 41007
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 1.5
$ws.Range("G12").Value = 1083
$ws.Range("H12").Value = 0.03847222222222222
$ws.Range("I12").Value = 6900
$ws.Range("J12").Value = "Vampiro"
$ws.Range("K12").Value = "Desafio"
$ws.Range("L12").Value = 46015

# D12 = C12*F12, following the same pattern as the shared formula used by
# the rows above (D2:D11).
$ws.Range("D12").Formula = "=C12*F12"

# --- Restore the selection to where the user last left off ---
$ws.Range("I13").Select()
